# edit.ps1 -- reproduces the commit that replaced the "Juan / DNI" sample
# data in Hoja1 with the "jorge / ID4" sample data, dropped the DNI/NIF/
# Nacionalidad/FechaNacimiento/PollingStation columns (F:I), removed the
# two mailto: hyperlinks that lived on C2/C4, and moved the active-cell
# selection to I8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two "mailto:juan@example.com" hyperlinks on C2/C4. This only
# removes the <hyperlinks> relationships; the cell values/styles are left
# alone (the "Hipervinculo" style (s=2) stays on C2/C4 per the target).
$ws.Hyperlinks.Delete()

# The new data set only needs columns A:E -- clear out the old F:I columns
# (Fecha nacimiento/Nacionalidad/DNI/NIF/pollingStation) entirely so the
# sheet's used range / dimension shrinks back down to A1:E4.
$ws.Range("F1:I4").ClearContents()

# Row 1: header row
$ws.Range("A1").Value = "Nombre"
$ws.Range("A2").Value = "jorge"
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("D1").Value = "ID4"
$ws.Range("B1").Value = "location"
$ws.Range("C1").Value = "email"
$ws.Range("E1").Value = "kind"

# Row 4 mirrors row 2 (row 3 stays empty, same as before the edit)
$ws.Range("A4").Value = "jorge"
$ws.Range("B4").Value = "18:13:14:12S"
$ws.Range("C4").Value = "jorge@email.es"
$ws.Range("D4").Value = "ID4"
$ws.Range("E2").Value = 1
$ws.Range("E4").Value = 1

# Move the selection the way the saved workbook shows it
$ws.Range("I8").Select()
